$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Environment Settings")

# Update the UAT URL cell (C3) to point to the new demo environment URL.
$ws.Range("C3").Value = "https://uatcrm-demo.ascentis.com.sg/AscentisCRM2/login.aspx"

# Update the active selection on the "Environment Settings" sheet from C10 to C12.
$ws.Activate()
$ws.Range("C12").Select()
